# Ajout gestion des filtres. Préparation sauvegarde des filtres.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for "Situation donateur" (column D) per row, cycling
# Grand donateur / Petit donateur / Moyen donateur / donateur nul
$situations = @{
    3  = "Petit donateur"
    4  = "Moyen donateur"
    5  = "donateur nul"
    7  = "Petit donateur"
    8  = "Moyen donateur"
    9  = "donateur nul"
    11 = "Petit donateur"
    12 = "Moyen donateur"
    13 = "donateur nul"
    15 = "Petit donateur"
    16 = "Moyen donateur"
    17 = "donateur nul"
    19 = "Petit donateur"
    20 = "Moyen donateur"
    21 = "donateur nul"
    23 = "Petit donateur"
    24 = "Moyen donateur"
    25 = "donateur nul"
    27 = "Petit donateur"
    28 = "Moyen donateur"
    29 = "donateur nul"
    31 = "Petit donateur"
    32 = "Moyen donateur"
    33 = "donateur nul"
}

foreach ($row in $situations.Keys) {
    $ws.Cells.Item($row, 4).Value = $situations[$row]
}

# Values for "COMPAGNIE" (column I) per row
$compagnies = @{
    2  = "Atos"
    3  = "Atos"
    4  = "Amadeus"
    5  = "Amadeus"
    6  = "Amadeus"
    7  = "Amadeus"
    8  = "Amadeus"
    9  = "Amadeus"
    10 = "Monext"
    11 = "Monext"
    12 = "Monext"
    13 = "Monext"
    14 = "Monext"
    15 = "Monext"
    16 = "Monext"
    17 = "Monext"
    18 = "Monext"
    19 = "Monext"
    20 = "Atos"
    21 = "Atos"
    22 = "Atos"
    23 = "Atos"
    24 = "Atos"
    25 = "Atos"
    26 = "Atos"
    27 = "Amadeus"
    28 = "Amadeus"
    29 = "Amadeus"
    30 = "Amadeus"
    31 = "Amadeus"
    32 = "Amadeus"
    33 = "Amadeus"
    34 = "Amadeus"
    35 = "Amadeus"
}

foreach ($row in $compagnies.Keys) {
    $ws.Cells.Item($row, 9).Value = $compagnies[$row]
}

# Update the active cell selection to I39 as recorded in the view state
$ws.Range("I39").Select()
